$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(" Oct 24 2020", " Abu Dhabi", "KKR won by 59 runs", "Delhi Capitals", "Kolkata Knight Riders", "Shimron Hetmyer ", "10", "5", "0", "1", "200.00"),
    @(" Oct 31 2020", " Dubai (DSC)", "Mumbai won by 9 wickets (with 34 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Shimron Hetmyer ", "11", "13", "1", "0", "84.61"),
    @(" Nov 10 2020", " Dubai (DSC)", "Mumbai won by 5 wickets (with 8 balls remaining)", "Delhi Capitals", "Mumbai Indians", "Shimron Hetmyer ", "5", "5", "1", "0", "100.00"),
    @(" Nov 8 2020", " Abu Dhabi", "Capitals won by 17 runs", "Delhi Capitals", "Sunrisers Hyderabad", "Shimron Hetmyer ", "42", "22", "4", "1", "190.90"),
    @(" Oct 9 2020", " Sharjah", "Capitals won by 46 runs", "Delhi Capitals", "Rajasthan Royals", "Shimron Hetmyer ", "45", "24", "1", "5", "187.50"),
    @(" Oct 5 2020", " Dubai (DSC)", "Capitals won by 59 runs", "Delhi Capitals", "Royal Challengers Bangalore", "Shimron Hetmyer ", "11", "7", "0", "1", "157.14"),
    @(" Sep 20 2020", " Dubai (DSC)", "Match tied (Capitals won the one-over eliminator)", "Delhi Capitals", "Kings XI Punjab", "Shimron Hetmyer ", "7", "13", "1", "0", "53.84"),
    @(" Oct 20 2020", " Dubai (DSC)", "Kings XI won by 5 wickets (with 6 balls remaining)", "Delhi Capitals", "Kings XI Punjab", "Shimron Hetmyer ", "10", "6", "0", "1", "166.66"),
    @(" Sep 29 2020", " Abu Dhabi", "Sunrisers won by 15 runs", "Delhi Capitals", "Sunrisers Hyderabad", "Shimron Hetmyer ", "21", "12", "0", "2", "175.00")
)

$lastRow = 1 + $data.Length
$numRange = $ws.Range("G2:K$lastRow")
# Force the numeric-looking values (runs/balls/4s/6s/sr) to be stored as TEXT,
# matching the source data (t="str" cells), not auto-converted to numbers.
$numRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

# Reset style back to the workbook's default "Normal" so no extra cell
# style/format is left behind beyond the value itself.
$numRange.Style = "Normal"
